# Update the "Förändrad" (changed) date column (C) for all data rows
# from serial date 45177 (2023-09-08) to 45178 (2023-09-09).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C470").Value = 45178
